# Insert a new data row at row 459 (shifts existing rows 459-523 down to 460-524,
# which also extends the sheet dimension from A1:T523 to A1:T524 automatically).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(459).Insert()

$ws.Range("A459").Value = 11
$ws.Range("B459").Value = "Vega Monumental Concepción"
$ws.Range("C459").Value = "Bíobío"
$ws.Range("D459").Value = 45142
$ws.Range("E459").Value = 8
$ws.Range("F459").Value = "Fruta"
$ws.Range("G459").Value = 100102
$ws.Range("H459").Value = "Cítricos"
$ws.Range("I459").Value = 100102005
$ws.Range("J459").Value = "Naranja"
$ws.Range("K459").Value = "Lane Late"
$ws.Range("L459").Value = "Primera"
$ws.Range("M459").Value = 100
$ws.Range("N459").Value = 7500
$ws.Range("O459").Value = 8000
$ws.Range("P459").Value = 7750
$ws.Range("Q459").Value = "$/caja 15 kilos empedrada"
$ws.Range("R459").Value = "Provincia de Melipilla"
$ws.Range("S459").Value = 517
$ws.Range("T459").Value = 15
